$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4   = -20.67159999999998
    6   = -22.65120000000001
    7   = -19.51909999999998
    16  = -21.6807
    20  = -19.8821
    28  = -22.0322
    29  = -21.42699999999997
    32  = -21.14369999999998
    40  = -20.04159999999999
    46  = -21.83060000000001
    51  = -21.96899999999999
    52  = -21.9888
    57  = -22.69830000000001
    59  = -22.1032
    62  = -22.0495
    66  = -21.5025
    73  = -20.02829999999999
    74  = -21.81139999999998
    92  = -21.42760000000002
    100 = -21.75229999999999
}

foreach ($row in $updates.Keys) {
    $ws.Range("A$row").Value = $updates[$row]
}
